# Daily attendance processing - 2026-01-22 06:09:43
# Reorders the comma-separated "Recorded By" names in column G for the
# session rows touched by this run (the order of the recorder names is
# swapped, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,
    24,26,29,32,33,36,37,38,39,40,41,43,44,45,46,
    47,48,50,52,55,58,59,62,63,64,65,66,67,69,70,
    71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,
    96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,
    135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = [string]$cell.Value2
    $parts = $current -split ', '
    $newValue = $parts[1] + ", " + $parts[0]
    $cell.Value = $newValue
}
